$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: the underlying molecule-name columns shifted by one position
# (values from B1:K1 each take on the meaning previously held one column to
# their "right" in the species list, wrapping "L" from A1 around to the end).
$ws.Range("B1").Value = "Formaldehyde"
$ws.Range("C1").Value = "14"
$ws.Range("D1").Value = "31-ol"
$ws.Range("E1").Value = "HOAc"
$ws.Range("F1").Value = "H2"
$ws.Range("G1").Value = "27-eth"
$ws.Range("H1").Value = "CO2"
$ws.Range("I1").Value = "CH4"
$ws.Range("J1").Value = "Water"
$ws.Range("K1").Value = "CO"

# Data rows 2-6: values re-aligned to the corrected column headers
# (HOAc uptake numbers recomputed given the Ni-uptake correction).
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 87879.26254673727
$ws.Range("F2").Value = 189242.8962731658
$ws.Range("H2").Value = 37447.79962493382
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 863.541570785506
$ws.Range("K2").Value = 6113.049791715604

$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 901.3547508711619
$ws.Range("F3").Value = 64662.57134536825
$ws.Range("H3").Value = 5247.788457850268
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 463.5429877008403
$ws.Range("K3").Value = 5028.771992537751

$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 355770.1455863371
$ws.Range("F4").Value = 1542796.678067265
$ws.Range("H4").Value = 282053.5960098541
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 973.8920375522388
$ws.Range("K4").Value = 7000.706858132517

$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 70150.2883145026
$ws.Range("F5").Value = 226436.697752606
$ws.Range("H5").Value = 70746.02542130341
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 379.9426308961204
$ws.Range("K5").Value = 5814.164763092483

$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 20014.14776842823
$ws.Range("F6").Value = 189890.7925692354
$ws.Range("H6").Value = 33415.81491505161
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 518.5141076504993
$ws.Range("K6").Value = 6619.894999744876
